# Update the cryptocurrency price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to a text number-format first so that price strings such as
# "69.413.63" or "179.86" are stored as literal text (matching the source feed)
# instead of being auto-coerced into numbers by Excel's Value setter.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.413.63'
$ws.Range("E2").Value = '  +2.11%  '

$ws.Range("D3").Value = '3.389.58'
$ws.Range("E3").Value = '  +1.42%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '587.87'
$ws.Range("E5").Value = '  +0.90%  '

$ws.Range("D6").Value = '179.86'
$ws.Range("E6").Value = '  +1.20%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '0.598'
$ws.Range("E8").Value = '  +1.26%  '

$ws.Range("E9").Value = '  +5.26%  '

$ws.Range("E10").Value = '  +1.30%  '

$ws.Range("D11").Value = '48.58'
$ws.Range("E11").Value = '  +2.59%  '

$ws.Range("E12").Value = '  +2.97%  '

$ws.Range("D13").Value = '679.93'
$ws.Range("E13").Value = '  -3.90%  '

$ws.Range("D14").Value = '8.62'
$ws.Range("E14").Value = '  +1.93%  '

$ws.Range("D15").Value = '3.933.94'
$ws.Range("E15").Value = '  +1.27%  '

$ws.Range("D16").Value = '69.501.00'
$ws.Range("E16").Value = '  +2.13%  '

$ws.Range("E17").Value = '  +1.70%  '

$ws.Range("D18").Value = '3.390.29'
$ws.Range("E18").Value = '  +1.40%  '

$ws.Range("E19").Value = '  +0.66%  '

$ws.Range("E20").Value = '  +1.74%  '

$ws.Range("D21").Value = '0.906'
$ws.Range("E21").Value = '  +0.88%  '

$ws.Range("D23").Value = '17.25'
$ws.Range("E23").Value = '  +0.75%  '

$ws.Range("D24").Value = '103.77'
$ws.Range("E24").Value = '  +3.38%  '

$ws.Range("E25").Value = '  +0.19%  '

$ws.Range("E26").Value = '  +0.57%  '

$ws.Range("D27").Value = '9.72'
$ws.Range("E27").Value = '  +0.76%  '

$ws.Range("D28").Value = '34.14'
$ws.Range("E28").Value = '  +2.82%  '

$ws.Range("D29").Value = '8.75'
$ws.Range("E29").Value = '  +1.53%  '

$ws.Range("D30").Value = '6.99'
$ws.Range("E30").Value = '  -1.12%  '

$ws.Range("E31").Value = '  +1.01%  '

$ws.Range("D32").Value = '559.96'
$ws.Range("E32").Value = '  -2.03%  '

$ws.Range("E33").Value = '  +0.57%  '

$ws.Range("D34").Value = '3.59'
$ws.Range("E34").Value = '  +4.55%  '

$ws.Range("D35").Value = '58.57'
$ws.Range("E35").Value = '  +0.98%  '

$ws.Range("E36").Value = '  +0.08%  '

$ws.Range("D37").Value = '3.693.85'
$ws.Range("E37").Value = '  -0.20%  '

$ws.Range("E38").Value = '  +4.75%  '

$ws.Range("D39").Value = '35.57'
$ws.Range("E39").Value = '  +2.42%  '

$ws.Range("E40").Value = '  +2.61%  '

$ws.Range("E41").Value = '  +1.44%  '

$ws.Range("D42").Value = '0.0₃0701'
$ws.Range("E42").Value = '  +3.27%  '

$ws.Range("E43").Value = '  +0.50%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0424'
$ws.Range("E44").Value = '  +3.74%  '

$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = '3.31'
$ws.Range("E45").Value = '  -2.02%  '

$ws.Range("E46").Value = '  -0.17%  '

$ws.Range("E47").Value = '  +0.83%  '

$ws.Range("D48").Value = '1.41'
$ws.Range("E48").Value = '  +5.07%  '

$ws.Range("E49").Value = '  -0.20%  '

$ws.Range("D50").Value = '133.19'
$ws.Range("E50").Value = '  +1.49%  '

$ws.Range("D51").Value = '2.62'
$ws.Range("E51").Value = '  +3.53%  '

# Restore the default (general) cell formatting on column D now that the text
# values are safely stored, so styling matches the original workbook.
$ws.Range("D2:D51").ClearFormats()
